$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header row labels so that the "_old" / "_new" suffixes used
#    to mark the two compared AHB format versions are replaced with the
#    concrete format version identifiers FV2310 (old/left side) and
#    FV2404 (new/right side). Column K ("diff") is untouched.
# ---------------------------------------------------------------------------
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2310")
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2404")
}

# ---------------------------------------------------------------------------
# 2. Turn the used range into an actual Excel Table ("Table1") spanning the
#    whole sheet (A1:U65), with a header row and an autofilter, so that the
#    columns carry the same (renamed) headers.
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U65")
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# ---------------------------------------------------------------------------
# 3. Freeze the header row (row 1) so it stays visible while scrolling.
# ---------------------------------------------------------------------------
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Renamed headers, added Table1 over A1:U65, froze header row."
